# Auto-generated Excel COM-interop script
# Applies numeric value updates to the Odin_Profits workbook sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the target diff.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 2035
$ws.Range("I5").Value = 1835.6666
$ws.Range("J5").Value = 2184.5
$ws.Range("K5").Value = 1835.6666
$ws.Range("L5").Value = 2184.5
$ws.Range("M5").Value = -1720.6666
$ws.Range("N5").Value = -2414.5
$ws.Range("H6").Value = 471
$ws.Range("I6").Value = 395.57144
$ws.Range("J6").Value = 999
$ws.Range("K6").Value = 1186.71432
$ws.Range("L6").Value = 2997
$ws.Range("M6").Value = -1074.71432
$ws.Range("N6").Value = -3221
$ws.Range("H8").Value = 42.875
$ws.Range("I8").Value = 41.857143
$ws.Range("J8").Value = 50
$ws.Range("K8").Value = 125.571429
$ws.Range("L8").Value = 150
$ws.Range("M8").Value = 13.42857100000001
$ws.Range("N8").Value = -428
$ws.Range("H40").Value = 11995
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 11995
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 11995
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -12345
$ws.Range("H52").Value = 5000
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 5000
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 15000
$ws.Range("N52").Value = -15320
$ws.Range("H88").Value = 3041.889
$ws.Range("I88").Value = 1860
$ws.Range("J88").Value = 3632.8333
$ws.Range("K88").Value = 1860
$ws.Range("L88").Value = 3632.8333
$ws.Range("M88").Value = -1454
$ws.Range("N88").Value = -4444.8333
$ws.Range("H91").Value = 3041.889
$ws.Range("I91").Value = 1860
$ws.Range("J91").Value = 3632.8333
$ws.Range("K91").Value = 1860
$ws.Range("L91").Value = 3632.8333
$ws.Range("M91").Value = -456
$ws.Range("N91").Value = -6440.8333
$ws.Range("H125").Value = 2668.8044
$ws.Range("I125").Value = 2377.5518
$ws.Range("J125").Value = 3165.647
$ws.Range("K125").Value = 21397.9662
$ws.Range("L125").Value = 28490.823
$ws.Range("M125").Value = -18937.9662
$ws.Range("N125").Value = -33410.823
$ws.Range("H132").Value = 224530.78
$ws.Range("I132").Value = 231478.3
$ws.Range("J132").Value = 30000
$ws.Range("K132").Value = 694434.8999999999
$ws.Range("L132").Value = 90000
$ws.Range("M132").Value = -691904.8999999999
$ws.Range("N132").Value = -95060
$ws.Range("H138").Value = 2759.8396
$ws.Range("I138").Value = 879.3261
$ws.Range("J138").Value = 5231.3716
$ws.Range("K138").Value = 2637.9783
$ws.Range("L138").Value = 15694.1148
$ws.Range("M138").Value = 2502.0217
$ws.Range("N138").Value = -25974.1148
$ws.Range("H141").Value = 1082.2759
$ws.Range("I141").Value = 1123.7307
$ws.Range("J141").Value = 723
$ws.Range("K141").Value = 3371.1921
$ws.Range("L141").Value = 2169
$ws.Range("M141").Value = 1808.8079
$ws.Range("N141").Value = -12529

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 359.2857
$ws.Range("I4").Value = 246.25
$ws.Range("J4").Value = 510
$ws.Range("K4").Value = 246.25
$ws.Range("L4").Value = 510
$ws.Range("M4").Value = -130.25
$ws.Range("N4").Value = -742
$ws.Range("H32").Value = 7940591.5
$ws.Range("I32").Value = 7464903.5
$ws.Range("J32").Value = 11127700
$ws.Range("K32").Value = 7464903.5
$ws.Range("L32").Value = 11127700
$ws.Range("M32").Value = -7464616.5
$ws.Range("N32").Value = -11128274
$ws.Range("H63").Value = 1730.4828
$ws.Range("I63").Value = 2252.7334
$ws.Range("J63").Value = 1170.9286
$ws.Range("K63").Value = 2252.7334
$ws.Range("L63").Value = 1170.9286
$ws.Range("M63").Value = -1566.7334
$ws.Range("N63").Value = -2542.9286
$ws.Range("H66").Value = 1730.4828
$ws.Range("I66").Value = 2252.7334
$ws.Range("J66").Value = 1170.9286
$ws.Range("K66").Value = 11263.667
$ws.Range("L66").Value = 5854.643
$ws.Range("M66").Value = -7831.667000000001
$ws.Range("N66").Value = -12718.643
$ws.Range("H74").Value = 4732.6523
$ws.Range("I74").Value = 4804.9375
$ws.Range("J74").Value = 4567.4287
$ws.Range("K74").Value = 4804.9375
$ws.Range("L74").Value = 4567.4287
$ws.Range("M74").Value = -3930.9375
$ws.Range("N74").Value = -6315.4287
$ws.Range("H77").Value = 4732.6523
$ws.Range("I77").Value = 4804.9375
$ws.Range("J77").Value = 4567.4287
$ws.Range("K77").Value = 24024.6875
$ws.Range("L77").Value = 22837.1435
$ws.Range("M77").Value = -19656.6875
$ws.Range("N77").Value = -31573.1435
$ws.Range("H132").Value = 498997.47
$ws.Range("I132").Value = 551313.2
$ws.Range("J132").Value = 126248
$ws.Range("K132").Value = 1653939.6
$ws.Range("L132").Value = 378744
$ws.Range("M132").Value = -1651409.6
$ws.Range("N132").Value = -383804

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 16176.294
$ws.Range("I86").Value = 14064.875
$ws.Range("J86").Value = 18053.111
$ws.Range("K86").Value = 14064.875
$ws.Range("L86").Value = 18053.111
$ws.Range("M86").Value = -12941.875
$ws.Range("N86").Value = -20299.111
$ws.Range("H89").Value = 16176.294
$ws.Range("I89").Value = 14064.875
$ws.Range("J89").Value = 18053.111
$ws.Range("K89").Value = 70324.375
$ws.Range("L89").Value = 90265.55500000001
$ws.Range("M89").Value = -64708.375
$ws.Range("N89").Value = -101497.555
$ws.Range("H134").Value = 762335.25
$ws.Range("I134").Value = 823510.3
$ws.Range("J134").Value = 15999.8
$ws.Range("K134").Value = 2470530.9
$ws.Range("L134").Value = 47999.39999999999
$ws.Range("M134").Value = -2467995.9
$ws.Range("N134").Value = -53069.39999999999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5449.3335
$ws.Range("I31").Value = 2837.5833
$ws.Range("J31").Value = 7538.7334
$ws.Range("K31").Value = 2837.5833
$ws.Range("L31").Value = 7538.7334
$ws.Range("M31").Value = -2542.5833
$ws.Range("N31").Value = -8128.7334
$ws.Range("H34").Value = 5449.3335
$ws.Range("I34").Value = 2837.5833
$ws.Range("J34").Value = 7538.7334
$ws.Range("K34").Value = 2837.5833
$ws.Range("L34").Value = 7538.7334
$ws.Range("M34").Value = -2635.5833
$ws.Range("N34").Value = -7942.7334
$ws.Range("H134").Value = 45460736
$ws.Range("I134").Value = 50003836
$ws.Range("J134").Value = 29750
$ws.Range("K134").Value = 150011508
$ws.Range("L134").Value = 89250
$ws.Range("M134").Value = -150008973
$ws.Range("N134").Value = -94320

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 2499.3333
$ws.Range("I25").Value = 1500
$ws.Range("J25").Value = 2999
$ws.Range("K25").Value = 4500
$ws.Range("L25").Value = 8997
$ws.Range("M25").Value = -4331
$ws.Range("N25").Value = -9335
$ws.Range("H26").Value = 249.36363
$ws.Range("I26").Value = 257.14285
$ws.Range("J26").Value = 235.75
$ws.Range("K26").Value = 771.4285500000001
$ws.Range("L26").Value = 707.25
$ws.Range("M26").Value = -483.4285500000001
$ws.Range("N26").Value = -1283.25
$ws.Range("H30").Value = 2499.3333
$ws.Range("I30").Value = 1500
$ws.Range("J30").Value = 2999
$ws.Range("K30").Value = 4500
$ws.Range("L30").Value = 8997
$ws.Range("M30").Value = -4398
$ws.Range("N30").Value = -9201
$ws.Range("H132").Value = 2923.4707
$ws.Range("I132").Value = 1952.4
$ws.Range("J132").Value = 3328.0833
$ws.Range("K132").Value = 17571.6
$ws.Range("L132").Value = 29952.7497
$ws.Range("M132").Value = -15041.6
$ws.Range("N132").Value = -35012.7497

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H49").Value = 24000
$ws.Range("I49").Value = 19000
$ws.Range("J49").Value = 26500
$ws.Range("K49").Value = 19000
$ws.Range("L49").Value = 26500
$ws.Range("M49").Value = -18816
$ws.Range("N49").Value = -26868
$ws.Range("H52").Value = 33000
$ws.Range("I52").Value = 33000
$ws.Range("J52").Value = 0
$ws.Range("K52").Value = 33000
$ws.Range("L52").Value = 0
$ws.Range("M52").Value = -32741
$ws.Range("H126").Value = 24008428
$ws.Range("I126").Value = 35719024
$ws.Range("J126").Value = 9104032
$ws.Range("K126").Value = 107157072
$ws.Range("L126").Value = 27312096
$ws.Range("M126").Value = -107154602
$ws.Range("N126").Value = -27317036
$ws.Range("H139").Value = 132720.75
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 132720.75
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 132720.75
$ws.Range("N139").Value = -143000.75

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2890.1755
$ws.Range("I132").Value = 3014.1304
$ws.Range("J132").Value = 2371.818
$ws.Range("K132").Value = 9042.3912
$ws.Range("L132").Value = 7115.454000000001
$ws.Range("M132").Value = -6512.3912
$ws.Range("N132").Value = -12175.454

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2382.25
$ws.Range("I126").Value = 2008.75
$ws.Range("J126").Value = 4249.75
$ws.Range("K126").Value = 6026.25
$ws.Range("L126").Value = 12749.25
$ws.Range("M126").Value = -3556.25
$ws.Range("N126").Value = -17689.25
$ws.Range("H132").Value = 6328.8486
$ws.Range("I132").Value = 4165
$ws.Range("J132").Value = 16066.167
$ws.Range("K132").Value = 12495
$ws.Range("L132").Value = 48198.501
$ws.Range("M132").Value = -9965
$ws.Range("N132").Value = -53258.501
